# Update countries & provincias Spain
# This script applies:
#  1) Country-name corrections in column A for rows whose ranking/order changed
#     (the underlying shared-string table is rebuilt by Excel in first-use order,
#     so updating the displayed names reproduces the reordering seen in the diff)
#  2) Updated statistic values (columns B-H) for the specific cells that changed
#  3) The "last updated" timestamp text in the title cell (A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Country name corrections (reordering) ---
    $ws.Cells.Item(20, 1).Value = "Canada"
    $ws.Cells.Item(21, 1).Value = "Australia"
    $ws.Cells.Item(22, 1).Value = "Brasil"
    $ws.Cells.Item(50, 1).Value = "Barein"
    $ws.Cells.Item(51, 1).Value = "Peru"
    $ws.Cells.Item(52, 1).Value = "Hong Kong"
    $ws.Cells.Item(53, 1).Value = "Estonia"
    $ws.Cells.Item(118, 1).Value = "Bolivia"
    $ws.Cells.Item(119, 1).Value = "Honduras"
    $ws.Cells.Item(127, 1).Value = "Guatemala"
    $ws.Cells.Item(128, 1).Value = "Guayana Francesa"
    $ws.Cells.Item(129, 1).Value = "Ruanda"
    $ws.Cells.Item(136, 1).Value = "Kirguistan"
    $ws.Cells.Item(137, 1).Value = "Kenia"
    $ws.Cells.Item(154, 1).Value = "Bahamas"
    $ws.Cells.Item(156, 1).Value = "Suazilandia"
    $ws.Cells.Item(159, 1).Value = "Zambia"
    $ws.Cells.Item(160, 1).Value = "Congo"
    $ws.Cells.Item(161, 1).Value = "San Bartolome"
    $ws.Cells.Item(162, 1).Value = "Republica de Africa Central"
    $ws.Cells.Item(163, 1).Value = "Fiyi"
    $ws.Cells.Item(164, 1).Value = "Santa Lucia"
    $ws.Cells.Item(165, 1).Value = "Liberia"
    $ws.Cells.Item(166, 1).Value = "El Salvador"
    $ws.Cells.Item(167, 1).Value = "Cabo Verde"
    $ws.Cells.Item(168, 1).Value = "Namibia"
    $ws.Cells.Item(171, 1).Value = "Benin"
    $ws.Cells.Item(172, 1).Value = "Niger"
    $ws.Cells.Item(173, 1).Value = "Mauritania"
    $ws.Cells.Item(174, 1).Value = "Butan"
    $ws.Cells.Item(176, 1).Value = "Birmania"
    $ws.Cells.Item(177, 1).Value = "Nicaragua"
    $ws.Cells.Item(179, 1).Value = "Nepal"
    $ws.Cells.Item(180, 1).Value = "Sudan"
    $ws.Cells.Item(182, 1).Value = "Montserrat"
    $ws.Cells.Item(184, 1).Value = "Republica de Yibuti"
    $ws.Cells.Item(185, 1).Value = "Islas Turcas y Caicos"
    $ws.Cells.Item(186, 1).Value = "Belice"
    $ws.Cells.Item(187, 1).Value = "Dominica"
    $ws.Cells.Item(188, 1).Value = "Republica del Chad"
    $ws.Cells.Item(189, 1).Value = "Eritrea"
    $ws.Cells.Item(190, 1).Value = "Timor Oriental"
    $ws.Cells.Item(191, 1).Value = "Papua Nueva Guinea"
    $ws.Cells.Item(192, 1).Value = "Uganda"
    $ws.Cells.Item(193, 1).Value = "Antigua y Barbuda"
    $ws.Cells.Item(194, 1).Value = "Mozambique"
    $ws.Cells.Item(195, 1).Value = "Santa Sede"
    $ws.Cells.Item(196, 1).Value = "Siria"
    $ws.Cells.Item(197, 1).Value = "San Vicente y las Granadinas"
    $ws.Cells.Item(198, 1).Value = "Granada"

# --- 2) Updated numeric statistics ---
    $ws.Range("B6").Value = 41526
    $ws.Range("C6").Value = 7980
    $ws.Range("E6").Value = 40838
    $ws.Range("G6").Value = 82
    $ws.Range("H6").Value = 501
    $ws.Range("E7").Value = 27527
    $ws.Range("G7").Value = 435
    $ws.Range("H7").Value = 2207
    $ws.Range("E10").Value = 13629
    $ws.Range("G10").Value = 186
    $ws.Range("H10").Value = 860
    $ws.Range("B13").Value = 6650
    $ws.Range("C13").Value = 967
    $ws.Range("E13").Value = 6180
    $ws.Range("B20").Value = 2020
    $ws.Range("C20").Value = 550
    $ws.Range("D20").Value = 18
    $ws.Range("E20").Value = 1979
    $ws.Range("F20").Value = 1
    $ws.Range("G20").Value = 3
    $ws.Range("H20").Value = 23
    $ws.Range("B21").Value = 1717
    $ws.Range("C21").Value = 108
    $ws.Range("D21").Value = 88
    $ws.Range("E21").Value = 1622
    $ws.Range("F21").Value = 11
    $ws.Range("H21").Value = 7
    $ws.Range("B22").Value = 1629
    $ws.Range("C22").Value = 83
    $ws.Range("D22").Value = 2
    $ws.Range("E22").Value = 1602
    $ws.Range("F22").Value = 18
    $ws.Range("G22").Value = 0
    $ws.Range("H22").Value = 25
    $ws.Range("B50").Value = 377
    $ws.Range("C50").Value = 43
    $ws.Range("D50").Value = 164
    $ws.Range("E50").Value = 211
    $ws.Range("F50").Value = 3
    $ws.Range("H50").Value = 2
    $ws.Range("B51").Value = 363
    $ws.Range("C51").Value = 0
    $ws.Range("D51").Value = 1
    $ws.Range("E51").Value = 357
    $ws.Range("F51").Value = 5
    $ws.Range("H51").Value = 5
    $ws.Range("B52").Value = 356
    $ws.Range("C52").Value = 38
    $ws.Range("D52").Value = 100
    $ws.Range("E52").Value = 252
    $ws.Range("H52").Value = 4
    $ws.Range("B53").Value = 352
    $ws.Range("C53").Value = 26
    $ws.Range("D53").Value = 4
    $ws.Range("E53").Value = 348
    $ws.Range("F53").Value = 4
    $ws.Range("H53").Value = 0
    $ws.Range("C118").Value = 3
    $ws.Range("C119").Value = 1
    $ws.Range("C127").Value = 1
    $ws.Range("D127").Value = 0
    $ws.Range("E127").Value = 19
    $ws.Range("H127").Value = 1
    $ws.Range("B128").Value = 20
    $ws.Range("C128").Value = 2
    $ws.Range("D128").Value = 6
    $ws.Range("E128").Value = 14
    $ws.Range("E129").Value = 19
    $ws.Range("H129").Value = 0
    $ws.Range("C136").Value = 2
    $ws.Range("C137").Value = 1
    $ws.Range("C161").Value = 0
    $ws.Range("C162").Value = 0
    $ws.Range("C163").Value = 1
    $ws.Range("C164").Value = 1
    $ws.Range("C171").Value = 0
    $ws.Range("C176").Value = 2
    $ws.Range("C179").Value = 1
    $ws.Range("D179").Value = 1
    $ws.Range("H179").Value = 0
    $ws.Range("C180").Value = 0
    $ws.Range("D180").Value = 0
    $ws.Range("H180").Value = 1
    $ws.Range("C185").Value = 1
    $ws.Range("C190").Value = 0

# --- 3) Updated "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 19:16"
